# edit.ps1 - PowerPoint COM-interop script reproducing the target diff.
#
# Summary of changes:
#  1. Slide 1 (title slide): title text becomes a two-line title
#     "二盐溶液溶解度曲面的" / "简捷算法"; subtitle author line becomes
#     "孙国铭 李泽健".
#  2. A new slide ("结论和展望") is inserted at position 10 (Title and
#     Content layout), pushing the closing "谢谢！" slide to position 11.
#  3. The "谢谢！" slide's text box is repositioned.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Title slide (slide 1): update title + subtitle text
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$titleShape = $s1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "二盐溶液溶解度曲面的`r简捷算法"

$subtitleShape = $s1.Shapes.Item(2)
$subtitleShape.TextFrame.TextRange.Paragraphs(1).Text = "孙国铭 李泽健"

# ---------------------------------------------------------------------
# 2. Insert new "结论和展望" slide at position 10 (Title and Content
#    layout = layout index 2), pushing the old last slide ("谢谢！") to
#    position 11.
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add(10, 2)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "结论和展望"

$newBody = $newSlide.Shapes.Item(2)
$bodyText = "计算效果很好`r这种方法使用简单的数据，简洁的模型，得到了与模拟软件精确计算相差无几的结果`r因此没有继续使用活度系数模型完善这部分工作`r如果后续继续完善的话，将活度系数模型引入，用活度代替现在工作中所有的浓度，可以得到更准确的结果。相应的，计算代价也会变大很多"
$newBody.TextFrame.TextRange.Text = $bodyText

# ---------------------------------------------------------------------
# 3. Reposition the "谢谢！" text box (now on the last slide).
# ---------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$thanksShape = $lastSlide.Shapes.Item(1)
$thanksShape.Left = 4705326 / 12700 + 0.00003
$thanksShape.Top = 2828835 / 12700 + 0.00003
